# TLWP-862 - Update opportunities report template
#
# The "Opportunities with providers" sheet's header row used to have:
#   G1 Provider contact name / H1 Provider contact email / I1 Provider contact telephone
#   J1 Secondary contact name / K1 Secondary contact email / L1 Secondary contact telephone
#
# It should become:
#   G1 Primary contact name / H1 Primary contact email / I1 Primary contact telephone
#   J1 Secondary contact name / K1 Secondary contact email / L1 Secondary contact telephone
#
# (the "Secondary contact ..." columns keep their position/content - only the
# former "Provider contact ..." columns are renamed to "Primary contact ...")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Opportunities with providers")

$ws.Range("G1").Value = "Primary contact name"
$ws.Range("H1").Value = "Primary contact email"
$ws.Range("I1").Value = "Primary contact telephone"
